$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Force used range / dimension recalculation
$dummy = $ws.UsedRange

# Copy number/date formatting from column E into the newly inserted column D
# (mirrors format of adjacent data so the new column matches the existing columns)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest-period figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 835800
$ws.Cells.Item(9, 4).Value = 591800
$ws.Cells.Item(10, 4).Value = 243900
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 11600
$ws.Cells.Item(15, 4).Value = 74600
$ws.Cells.Item(17, 4).Value = 761400
$ws.Cells.Item(18, 4).Value = 74400
$ws.Cells.Item(20, 4).Value = 49000
$ws.Cells.Item(21, 4).Value = 198000
$ws.Cells.Item(22, 4).Value = 31700
$ws.Cells.Item(23, 4).Value = 91700
$ws.Cells.Item(24, 4).Value = 8400
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 83300
$ws.Cells.Item(27, 4).Value = 58100
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -49000
$ws.Cells.Item(33, 4).Value = 58100
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 58100
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 144200
$ws.Cells.Item(42, 4).Value = 30300
$ws.Cells.Item(43, 4).Value = 210700
$ws.Cells.Item(44, 4).Value = 4500
$ws.Cells.Item(45, 4).Value = 8400
$ws.Cells.Item(46, 4).Value = 398100
$ws.Cells.Item(47, 4).Value = 161200
$ws.Cells.Item(48, 4).Value = 846500
$ws.Cells.Item(49, 4).Value = 57300
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 7900
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 1471000
$ws.Cells.Item(57, 4).Value = 59600
$ws.Cells.Item(58, 4).Value = 8500
$ws.Cells.Item(59, 4).Value = 55700
$ws.Cells.Item(60, 4).Value = 123800
$ws.Cells.Item(61, 4).Value = 346100
$ws.Cells.Item(62, 4).Value = 147300
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 766900
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 474800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 704200
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = 58100
$ws.Cells.Item(83, 4).Value = 74600
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 49400
$ws.Cells.Item(91, 4).Value = -50300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = 80500
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = -224800
$ws.Cells.Item(101, 4).Value = -100
$ws.Cells.Item(102, 4).Value = -95000

# A handful of prior-period figures were also restated alongside the new column
$ws.Cells.Item(8, 5).Value = 650800
$ws.Cells.Item(8, 6).Value = 524200
$ws.Cells.Item(9, 5).Value = 433800
$ws.Cells.Item(9, 6).Value = 359000
$ws.Cells.Item(17, 5).Value = 601200
$ws.Cells.Item(17, 6).Value = 522200
$ws.Cells.Item(94, 5).Value = 118600
$ws.Cells.Item(94, 6).Value = -109700
$ws.Cells.Item(102, 5).Value = -16700
$ws.Cells.Item(102, 6).Value = -98300

# Match column D width to the neighboring data columns
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth
